# Provider Dupe Logic Added
# Duplicate Providers can cause issues, so add a duplicate provider row
# and update the file output type of an existing provider, so the dupe
# detection logic has something to catch.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Change 1: update FILE_OUTPUT_TYPE for the existing GS23FBA011 / row 2
#     provider from "csv" to "xlsx"
$ws.Range("E2").Value = "xlsx"

# --- Change 2: add a new row 4 that duplicates the GS23FBA011 provider
#     (same PROVIDER_ID/IDENTIFIER/EMAIL as row 2) but keeps the original
#     "csv" output type, creating a duplicate provider entry.

# First copy the formatting (font/border/alignment) of an existing data
# row onto the new row so it matches the rest of the table (style index
# used by rows 2 and 3).
$ws.Range("A3:E3").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)

# Fill in the values for the new row.
$ws.Range("A4").Value = "001"
$ws.Range("B4").Value = "MRO"
$ws.Range("C4").Value = "GS23FBA011"
$ws.Range("D4").Value = "david.larrimore@gmail.com"
$ws.Range("E4").Value = "csv"

# Column A ("001") looks numeric, so a plain .Value assignment gets
# auto-coerced to the number 1, losing the leading zeros. Fix this by
# copying both the value and the format from an existing cell (A2) that
# already stores "001" correctly as text, which keeps it as text without
# disturbing the row's style.
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4104)
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
